$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "rana.abozaid@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("L7").Value = 23
$ws.Range("L8").Value = 96
$ws.Range("G10").Value = "arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg"
$ws.Range("G15").Value = "hana.amr@med.asu.edu.eg, nancy.abdelshafy@med.asu.edu.eg"
$ws.Range("G18").Value = "shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G19").Value = "Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Range("G20").Value = "aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Range("P21").Value = 3
$ws.Range("Q21").Value = 12
$ws.Range("P22").Value = 4
$ws.Range("Q22").Value = 12
$ws.Range("G24").Value = "rana.abozaid@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G32").Value = "arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg"
$ws.Range("G37").Value = "hana.amr@med.asu.edu.eg, nancy.abdelshafy@med.asu.edu.eg"
$ws.Range("G40").Value = "shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G41").Value = "Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Range("G42").Value = "aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Range("G46").Value = "nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G52").Value = "mariam.noureldin@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("G54").Value = "Eman.m.abosakaya@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, merna.said@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, maimustafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg"
$ws.Range("G58").Value = "afaf.abdallah@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg"
$ws.Range("G63").Value = "wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G68").Value = "nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G74").Value = "mariam.noureldin@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("G76").Value = "Eman.m.abosakaya@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, merna.said@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, maimustafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg"
$ws.Range("G80").Value = "afaf.abdallah@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg"
$ws.Range("G85").Value = "wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G96").Value = "Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G98").Value = "Eman.m.abosakaya@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, merna.said@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, maimustafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg"
$ws.Range("G106").Value = "neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg"
$ws.Range("G107").Value = "aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Range("G108").Value = "aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Range("G118").Value = "Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G120").Value = "Eman.m.abosakaya@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, merna.said@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, maimustafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg"
$ws.Range("G128").Value = "neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg"
$ws.Range("G129").Value = "aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Range("G130").Value = "aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Range("G134").Value = "Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G142").Value = "marwa_mustafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, basma.hamed@med.asu.edu.eg"
$ws.Range("G150").Value = "Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Range("G156").Value = "manar.montaser@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G164").Value = "marwa_mustafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, basma.hamed@med.asu.edu.eg"
$ws.Range("G172").Value = "Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"

# Row 152 and 169: style change (Pending -> Not Recorded style) + status text update
$ws.Range("A3:I3").Copy()
$ws.Range("A152:I152").PasteSpecial(-4122)
$ws.Range("A3:I3").Copy()
$ws.Range("A169:I169").PasteSpecial(-4122)
$ws.Range("I152").Value = "Not Recorded"
$ws.Range("I169").Value = "Not Recorded"
